# Aggiornamento dati Carpi al 23 agosto 2021
# Appends daily rows (date, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.)
# for 2021-08-10 .. 2021-08-23 (Excel serial dates 44418..44431) to the bottom of Sheet1,
# extending the used range from A1:D343 to A1:D357.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 343
$firstNewRow = $lastRow + 1

# Carry the existing row's formatting (date style on col A, etc.) down onto the new rows.
$ws.Range("A343:D343").Copy() | Out-Null
$ws.Range("A344:D357").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$data = @(
  @(44418, 3, 68, 93.96288466055908),
  @(44419, 5, 71, 98.10830604264255),
  @(44420, 17, 78, 107.7809559341707),
  @(44421, 5, 70, 96.7264989152814),
  @(44422, 25, 77, 106.3991488068095),
  @(44423, 2, 74, 102.2537274247261),
  @(44424, 15, 72, 99.49011317000372),
  @(44425, 6, 75, 103.6355345520872),
  @(44426, 8, 78, 107.7809559341707),
  @(44427, 11, 72, 99.49011317000372),
  @(44428, 5, 72, 99.49011317000372),
  @(44429, 16, 63, 87.05384902375327),
  @(44430, 6, 67, 92.58107753319791),
  @(44431, 12, 64, 88.43565615111443)
)

$r = $firstNewRow
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}
